$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.053505811508217
$ws.Range("C2").Value = 0.1992382023926211
$ws.Range("D2").Value = 0.2241980064248708
$ws.Range("F2").Value = 1.089944437231786
$ws.Range("G2").Value = 0.5389779156351295
$ws.Range("H2").Value = 0.6749900868873766
$ws.Range("I2").Value = 0.6393562858807513
$ws.Range("J2").Value = 0.1780621131439286
$ws.Range("L2").Value = 0.3876720284449675
$ws.Range("N2").Value = 1.175515594943626
$ws.Range("O2").Value = 2.398732720235103

$ws.Range("B3").Value = 0.9620794073420598
$ws.Range("C3").Value = 0.1899420474211269
$ws.Range("D3").Value = 0.2234080634790772
$ws.Range("F3").Value = 1.091953429321521
$ws.Range("G3").Value = 0.5388811009410333
$ws.Range("H3").Value = 0.6786162384583392
$ws.Range("I3").Value = 0.645743417936977
$ws.Range("J3").Value = 0.1796287380230268
$ws.Range("L3").Value = 0.3794413889013413
$ws.Range("N3").Value = 1.178167065373636
$ws.Range("O3").Value = 2.405854679637031

$ws.Range("B4").Value = 0.9059858800534357
$ws.Range("C4").Value = 0.1841997647119626
$ws.Range("D4").Value = 0.2229982347096993
$ws.Range("F4").Value = 1.093771950556572
$ws.Range("G4").Value = 0.539189518556924
$ws.Range("H4").Value = 0.6811409513148021
$ws.Range("I4").Value = 0.6500110477806906
$ws.Range("J4").Value = 0.1806545973539446
$ws.Range("L4").Value = 0.3745257404509488
$ws.Range("N4").Value = 1.180209076753961
$ws.Range("O4").Value = 2.411622324664336

$ws.Range("B5").Value = 0.883139936694846
$ws.Range("C5").Value = 0.1818512720258383
$ws.Range("D5").Value = 0.222850208258798
$ws.Range("F5").Value = 1.094660225333982
$ws.Range("G5").Value = 0.5394076752139298
$ws.Range("H5").Value = 0.6822448567885857
$ws.Range("I5").Value = 0.6518371105565954
$ws.Range("J5").Value = 0.1810887378344148
$ws.Range("L5").Value = 0.3725574983989475
$ws.Range("N5").Value = 1.181145573082816
$ws.Range("O5").Value = 2.414323483140549

$ws.Range("B6").Value = 0.8793472035891625
$ws.Range("C6").Value = 0.1814608007593961
$ws.Range("D6").Value = 0.2228267771853254
$ws.Range("F6").Value = 1.094816617242806
$ws.Range("G6").Value = 0.539449484337851
$ws.Range("H6").Value = 0.6824326952017969
$ws.Range("I6").Value = 0.6521455792595177
$ws.Range("J6").Value = 0.181161799035392
$ws.Range("L6").Value = 0.3722327892547668
$ws.Range("N6").Value = 1.181307388311097
$ws.Range("O6").Value = 2.414793198014308

$ws.Range("B7").Value = 0.9056777180387598
$ws.Range("C7").Value = 0.1841681261496149
$ws.Range("D7").Value = 0.2229961614192035
$ws.Range("F7").Value = 1.093783333934596
$ws.Range("G7").Value = 0.5391920863076436
$ws.Range("H7").Value = 0.6811555349448497
$ws.Range("I7").Value = 0.6500353225832463
$ws.Range("J7").Value = 0.1806603871381363
$ws.Range("L7").Value = 0.3744990543440707
$ws.Range("N7").Value = 1.18022128377698
$ws.Range("O7").Value = 2.411657332967309

$ws.Range("B8").Value = 1.021974244160901
$ws.Range("C8").Value = 0.1960401377800594
$ws.Range("D8").Value = 0.2239100693617502
$ws.Range("F8").Value = 1.090515787356956
$ws.Range("G8").Value = 0.5388681656978207
$ws.Range("H8").Value = 0.676178525098031
$ws.Range("I8").Value = 0.6414867696361419
$ws.Range("J8").Value = 0.1785890226061859
$ws.Range("L8").Value = 0.3848055936836943
$ws.Range("N8").Value = 1.176344048110984
$ws.Range("O8").Value = 2.40089893800932

$ws.Range("B9").Value = 1.250297513892406
$ws.Range("C9").Value = 0.2190409988657223
$ws.Range("D9").Value = 0.2262959867447307
$ws.Range("F9").Value = 1.08874600079843
$ws.Range("G9").Value = 0.5411538878790907
$ws.Range("H9").Value = 0.6687822527688922
$ws.Range("I9").Value = 0.6274682467506771
$ws.Range("J9").Value = 0.1750337668216897
$ws.Range("L9").Value = 0.4061030906527208
$ws.Range("N9").Value = 1.172015239865729
$ws.Range("O9").Value = 2.390867283439661

$ws.Range("B10").Value = 1.418128021485416
$ws.Range("C10").Value = 0.2357609244943148
$ws.Range("D10").Value = 0.2284073367851107
$ws.Range("F10").Value = 1.090269667563291
$ws.Range("G10").Value = 0.5446180216760297
$ws.Range("H10").Value = 0.6647857972259885
$ws.Range("I10").Value = 0.6188431986386362
$ws.Range("J10").Value = 0.1727296183461657
$ws.Range("L10").Value = 0.422403091685922
$ws.Range("N10").Value = 1.170818240456626
$ws.Range("O10").Value = 2.390245045757155

$ws.Range("B11").Value = 1.494479993581365
$ws.Range("C11").Value = 0.2433268146108105
$ws.Range("D11").Value = 0.2294449434378265
$ws.Range("F11").Value = 1.091575387524927
$ws.Range("G11").Value = 0.5465824545558604
$ws.Range("H11").Value = 0.6632791878578672
$ws.Range("I11").Value = 0.6152832035624769
$ws.Range("J11").Value = 0.1717480480991824
$ws.Range("L11").Value = 0.4299582029079687
$ws.Range("N11").Value = 1.170701804173945
$ws.Range("O11").Value = 2.39142793466965

$ws.Range("B12").Value = 1.523391501343951
$ws.Range("C12").Value = 0.2461858860903021
$ws.Range("D12").Value = 0.2298488827573522
$ws.Range("F12").Value = 1.092157836699982
$ws.Range("G12").Value = 0.5473822633397276
$ws.Range("H12").Value = 0.6627533925965849
$ws.Range("I12").Value = 0.6139874282938571
$ws.Range("J12").Value = 0.1713859149384636
$ws.Range("L12").Value = 0.432839079159578
$ws.Range("N12").Value = 1.1707190493007
$ws.Range("O12").Value = 2.392086655667555

$ws.Range("B13").Value = 1.517164982489476
$ws.Range("C13").Value = 0.2455704024880561
$ws.Range("D13").Value = 0.2297613982319291
$ws.Range("F13").Value = 1.092028483398451
$ws.Range("G13").Value = 0.5472075224475503
$ws.Range("H13").Value = 0.6628646437125099
$ws.Range("I13").Value = 0.6142641692918218
$ws.Range("J13").Value = 0.1714634815032099
$ws.Range("L13").Value = 0.432217748728192
$ws.Range("N13").Value = 1.170712610748211
$ws.Range("O13").Value = 2.391935413793249

$ws.Range("B14").Value = 1.496858599055088
$ws.Range("C14").Value = 0.2435621528845218
$ws.Range("D14").Value = 0.2294779553244837
$ws.Range("F14").Value = 1.091621542806081
$ws.Range("G14").Value = 0.5466471343296604
$ws.Range("H14").Value = 0.6632350343172675
$ws.Range("I14").Value = 0.61517555063012
$ws.Range("J14").Value = 0.1717180635261535
$ws.Range("L14").Value = 0.4301948164591778
$ws.Range("N14").Value = 1.170701994906366
$ws.Range("O14").Value = 2.391477903661041

$ws.Range("B15").Value = 1.484420119245158
$ws.Range("C15").Value = 0.2423312589511397
$ws.Range("D15").Value = 0.2293057713529691
$ws.Range("F15").Value = 1.091383737649508
$ws.Range("G15").Value = 0.5463111642085465
$ws.Range("H15").Value = 0.6634677321796687
$ws.Range("I15").Value = 0.6157406123405167
$ws.Range("J15").Value = 0.1718752478529062
$ws.Range("L15").Value = 0.4289582991750223
$ws.Range("N15").Value = 1.17070347379989
$ws.Range("O15").Value = 2.391225115641618

$ws.Range("B16").Value = 1.413138166435544
$ws.Range("C16").Value = 0.2352656518053209
$ws.Range("D16").Value = 0.228341072834283
$ws.Range("F16").Value = 1.090196655095525
$ws.Range("G16").Value = 0.5444974643368852
$ws.Range("H16").Value = 0.6648905192332109
$ws.Range("I16").Value = 0.6190831718897911
$ws.Range("J16").Value = 0.1727951054946413
$ws.Range("L16").Value = 0.4219121476618
$ws.Range("N16").Value = 1.17083444500777
$ws.Range("O16").Value = 2.390197238603463

$ws.Range("B17").Value = 1.369408825965024
$ws.Range("C17").Value = 0.2309207190947689
$ws.Range("D17").Value = 0.2277689678439287
$ws.Range("F17").Value = 1.089625255886546
$ws.Range("G17").Value = 0.5434843743277042
$ws.Range("H17").Value = 0.665843071390924
$ws.Range("I17").Value = 0.6212268710579849
$ws.Range("J17").Value = 0.1733764582690398
$ws.Range("L17").Value = 0.4176252927579611
$ws.Range("N17").Value = 1.17102427573974
$ws.Range("O17").Value = 2.38994212416145

$ws.Range("B18").Value = 1.344257461956772
$ws.Range("C18").Value = 0.2284178675163844
$ws.Range("D18").Value = 0.2274471735973123
$ws.Range("F18").Value = 1.089354268047082
$ws.Range("G18").Value = 0.5429382443421815
$ws.Range("H18").Value = 0.6664202676899009
$ws.Range("I18").Value = 0.6224940901711911
$ws.Range("J18").Value = 0.1737171061987013
$ws.Range("L18").Value = 0.4151728136876187
$ws.Range("N18").Value = 1.171173767658033
$ws.Range("O18").Value = 2.389933372123437

$ws.Range("B19").Value = 1.335741807176021
$ws.Range("C19").Value = 0.2275698051611243
$ws.Range("D19").Value = 0.2273394696560018
$ws.Range("F19").Value = 1.089272423090897
$ws.Range("G19").Value = 0.5427596141516347
$ws.Range("H19").Value = 0.6666207326957618
$ws.Range("I19").Value = 0.6229290250476751
$ws.Range("J19").Value = 0.1738335207981265
$ws.Range("L19").Value = 0.4143447218027489
$ws.Range("N19").Value = 1.171231312505412
$ws.Range("O19").Value = 2.389954107960108

$ws.Range("B20").Value = 1.374063842086343
$ws.Range("C20").Value = 0.2313836353369538
$ws.Range("D20").Value = 0.2278291179818837
$ws.Range("F20").Value = 1.089680114749044
$ws.Range("G20").Value = 0.5435884340452759
$ws.Range("H20").Value = 0.665738637139782
$ws.Range("I20").Value = 0.6209951287212405
$ws.Range("J20").Value = 0.1733139235127474
$ws.Range("L20").Value = 0.4180802706620312
$ws.Range("N20").Value = 1.170999897997703
$ws.Range("O20").Value = 2.389955000899761

$ws.Range("B21").Value = 1.502823128444788
$ws.Range("C21").Value = 0.2441521880233211
$ws.Range("D21").Value = 0.2295609108726921
$ws.Range("F21").Value = 1.091738683507884
$ws.Range("G21").Value = 0.5468102158952206
$ws.Range("H21").Value = 0.6631250282219412
$ws.Range("I21").Value = 0.6149064356780478
$ws.Range("J21").Value = 0.1716430270751523
$ws.Range("L21").Value = 0.4307884619649229
$ws.Range("N21").Value = 1.170703450107908
$ws.Range("O21").Value = 2.391606564922455

$ws.Range("B22").Value = 1.586966081261778
$ws.Range("C22").Value = 0.2524623020922263
$ws.Range("D22").Value = 0.230756938010316
$ws.Range("F22").Value = 1.093596974226543
$ws.Range("G22").Value = 0.5492418041718707
$ws.Range("H22").Value = 0.6616775503484433
$ws.Range("I22").Value = 0.6112320857505225
$ws.Range("J22").Value = 0.1706067538670144
$ws.Range("L22").Value = 0.4392099971438768
$ws.Range("N22").Value = 1.170867123781321
$ws.Range("O22").Value = 2.393914574190603

$ws.Range("B23").Value = 1.54205887506987
$ws.Range("C23").Value = 0.2480302954410263
$ws.Range("D23").Value = 0.2301127455438774
$ws.Range("F23").Value = 1.092558269178269
$ws.Range("G23").Value = 0.5479141797829357
$ws.Range("H23").Value = 0.6624262632639955
$ws.Range("I23").Value = 0.6131652399036653
$ws.Range("J23").Value = 0.1711547343470645
$ws.Range("L23").Value = 0.4347047292364579
$ws.Range("N23").Value = 1.170747136091421
$ws.Range("O23").Value = 2.392570333430342

$ws.Range("B24").Value = 1.37195934272944
$ws.Range("C24").Value = 0.2311743661269361
$ws.Range("D24").Value = 0.2278019019322386
$ws.Range("F24").Value = 1.089655133889764
$ws.Range("G24").Value = 0.5435412755509219
$ws.Range("H24").Value = 0.6657857597785579
$ws.Range("I24").Value = 0.6210997911038199
$ws.Range("J24").Value = 0.1733421754760869
$ws.Range("L24").Value = 0.4178745374732387
$ws.Range("N24").Value = 1.171010793453149
$ws.Range("O24").Value = 2.389948749725647

$ws.Range("B25").Value = 1.188509898674511
$ws.Range("C25").Value = 0.2128494779753112
$ws.Range("D25").Value = 0.2255873350875461
$ws.Range("F25").Value = 1.088728611653735
$ws.Range("G25").Value = 0.540222401799241
$ws.Range("H25").Value = 0.6705304438041964
$ws.Range("I25").Value = 0.6309666527907716
$ws.Range("J25").Value = 0.1759414163872215
$ws.Range("L25").Value = 0.4002262357602575
$ws.Range("N25").Value = 1.172837087445615
$ws.Range("O25").Value = 2.392396237006039
